# Update column F (dSF) values for the specified rows to reflect
# repulled data / recalculated mean as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 6
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -3
